# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice* / LevePrice* / LeveProfit*) across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 60928.31
$ws.Range("I98").Value = 1026.4546
$ws.Range("J98").Value = 390388.5
$ws.Range("K98").Value = 1026.4546
$ws.Range("L98").Value = 390388.5
$ws.Range("M98").Value = 471.5454
$ws.Range("N98").Value = -393384.5

# Row 116
$ws.Range("H116").Value = 2233.4473
$ws.Range("I116").Value = 1508.5416
$ws.Range("J116").Value = 3476.1428
$ws.Range("K116").Value = 1508.5416
$ws.Range("L116").Value = 3476.1428
$ws.Range("M116").Value = 1933.4584
$ws.Range("N116").Value = -10360.1428

# Row 122
$ws.Range("H122").Value = 60928.31
$ws.Range("I122").Value = 1026.4546
$ws.Range("J122").Value = 390388.5
$ws.Range("K122").Value = 3079.3638
$ws.Range("L122").Value = 1171165.5
$ws.Range("M122").Value = -629.3638000000001
$ws.Range("N122").Value = -1176065.5

# Row 132
$ws.Range("H132").Value = 35983.86
$ws.Range("I132").Value = 23760.979
$ws.Range("J132").Value = 93431.39999999999
$ws.Range("K132").Value = 71282.93700000001
$ws.Range("L132").Value = 280294.2
$ws.Range("M132").Value = -68752.93700000001
$ws.Range("N132").Value = -285354.2

# Row 137
$ws.Range("H137").Value = 1795815.5
$ws.Range("I137").Value = 6421418.5
$ws.Range("J137").Value = 5259.4517
$ws.Range("K137").Value = 19264255.5
$ws.Range("L137").Value = 15778.3551
$ws.Range("M137").Value = -19261705.5
$ws.Range("N137").Value = -20878.3551

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12171.412
$ws.Range("I32").Value = 12710.579
$ws.Range("J32").Value = 10595.385
$ws.Range("K32").Value = 12710.579
$ws.Range("L32").Value = 10595.385
$ws.Range("M32").Value = -12423.579
$ws.Range("N32").Value = -11169.385

# Row 45
$ws.Range("H45").Value = 47620876
$ws.Range("J45").Value = 2158
$ws.Range("L45").Value = 2158
$ws.Range("N45").Value = -2912

# Row 62
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248

# Row 65
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240

# Row 122
$ws.Range("H122").Value = 1858.8572
$ws.Range("I122").Value = 1858.8572
$ws.Range("K122").Value = 5576.571599999999
$ws.Range("M122").Value = -3126.571599999999

# Row 132
$ws.Range("H132").Value = 12198016
$ws.Range("I132").Value = 20835798
$ws.Range("J132").Value = 3500.2354
$ws.Range("K132").Value = 62507394
$ws.Range("L132").Value = 10500.7062
$ws.Range("M132").Value = -62504864
$ws.Range("N132").Value = -15560.7062

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2364.4048
$ws.Range("I134").Value = 2268.9429
$ws.Range("K134").Value = 6806.8287
$ws.Range("M134").Value = -4271.8287

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1420.22
$ws.Range("I31").Value = 608.05457
$ws.Range("J31").Value = 2412.8667
$ws.Range("K31").Value = 608.05457
$ws.Range("L31").Value = 2412.8667
$ws.Range("M31").Value = -313.05457
$ws.Range("N31").Value = -3002.8667

# Row 34
$ws.Range("H34").Value = 1420.22
$ws.Range("I34").Value = 608.05457
$ws.Range("J34").Value = 2412.8667
$ws.Range("K34").Value = 608.05457
$ws.Range("L34").Value = 2412.8667
$ws.Range("M34").Value = -406.05457
$ws.Range("N34").Value = -2816.8667

# Row 134
$ws.Range("H134").Value = 803702.5600000001
$ws.Range("I134").Value = 564491.8
$ws.Range("J134").Value = 1401729.5
$ws.Range("K134").Value = 1693475.4
$ws.Range("L134").Value = 4205188.5
$ws.Range("M134").Value = -1690940.4
$ws.Range("N134").Value = -4210258.5

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 902.5833
$ws.Range("I8").Value = 902.5833
$ws.Range("K8").Value = 2707.7499
$ws.Range("M8").Value = -2568.7499

# Row 12
$ws.Range("H12").Value = 138.92857
$ws.Range("I12").Value = 162.6923
$ws.Range("J12").Value = 118.333336
$ws.Range("K12").Value = 488.0769
$ws.Range("L12").Value = 355.000008
$ws.Range("M12").Value = -315.0769
$ws.Range("N12").Value = -701.000008

# Row 23
$ws.Range("H23").Value = 699.0476
$ws.Range("I23").Value = 749.6
$ws.Range("J23").Value = 683.25
$ws.Range("K23").Value = 2248.8
$ws.Range("L23").Value = 2049.75
$ws.Range("M23").Value = -2013.8
$ws.Range("N23").Value = -2519.75

# Row 68
$ws.Range("H68").Value = 1117.2916
$ws.Range("J68").Value = 1241.8182
$ws.Range("L68").Value = 3725.4546
$ws.Range("N68").Value = -5347.4546

# Row 71
$ws.Range("H71").Value = 1117.2916
$ws.Range("J71").Value = 1241.8182
$ws.Range("L71").Value = 11176.3638
$ws.Range("N71").Value = -19288.3638

# Row 107
$ws.Range("H107").Value = 3873.1694
$ws.Range("I107").Value = 2899.9473
$ws.Range("J107").Value = 5634.2383
$ws.Range("K107").Value = 8699.841899999999
$ws.Range("L107").Value = 16902.7149
$ws.Range("M107").Value = -6779.841899999999
$ws.Range("N107").Value = -20742.7149

# Row 131
$ws.Range("H131").Value = 3566.4443
$ws.Range("I131").Value = 11577.667
$ws.Range("J131").Value = 1563.6389
$ws.Range("K131").Value = 34733.001
$ws.Range("L131").Value = 4690.9167
$ws.Range("M131").Value = -29693.001
$ws.Range("N131").Value = -14770.9167

# Row 140
$ws.Range("H140").Value = 57608.777
$ws.Range("I140").Value = 91756.37
$ws.Range("K140").Value = 275269.11
$ws.Range("M140").Value = -270089.11

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4135.385
$ws.Range("I40").Value = 3534.0557
$ws.Range("J40").Value = 5488.375
$ws.Range("K40").Value = 3534.0557
$ws.Range("L40").Value = 5488.375
$ws.Range("M40").Value = -3398.0557
$ws.Range("N40").Value = -5760.375

# Row 122
$ws.Range("H122").Value = 2365.6667
$ws.Range("I122").Value = 2417.0908
$ws.Range("K122").Value = 7251.2724
$ws.Range("M122").Value = -4801.2724

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 50050040
$ws.Range("I54").Value = 100000000
$ws.Range("K54").Value = 100000000
$ws.Range("M54").Value = -99999480

# Row 62
$ws.Range("H62").Value = 2999.9285
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -4248

# Row 65
$ws.Range("H65").Value = 2999.9285
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -21240

# Row 107
$ws.Range("H107").Value = 3846978.8
$ws.Range("J107").Value = 12500977
$ws.Range("L107").Value = 37502931
$ws.Range("N107").Value = -37506771
